$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "68.550.08"
$ws.Cells.Item(2, 5).Value = "  -1.37%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.454.97"
$ws.Cells.Item(3, 5).Value = "  -1.93%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "563.64"
$ws.Cells.Item(5, 5).Value = "  -2.02%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "162.99"
$ws.Cells.Item(6, 5).Value = "  -2.11%  "
$ws.Cells.Item(7, 5).Value = "  +0.03%  "
$ws.Cells.Item(8, 5).Value = "  -1.77%  "
$ws.Cells.Item(9, 2).Value = "LidoStakedEther"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "2.454.46"
$ws.Cells.Item(9, 5).Value = "  -1.88%  "
$ws.Cells.Item(10, 2).Value = "Dogecoin"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.152"
$ws.Cells.Item(10, 5).Value = "  -5.74%  "
$ws.Cells.Item(11, 2).Value = "TRON"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.164"
$ws.Cells.Item(11, 5).Value = "  -1.74%  "
$ws.Cells.Item(12, 2).Value = "Cardano"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.344"
$ws.Cells.Item(12, 5).Value = "  -4.08%  "
$ws.Cells.Item(13, 2).Value = "Toncoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.81"
$ws.Cells.Item(13, 5).Value = "  -2.12%  "
$ws.Cells.Item(14, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "2.905.57"
$ws.Cells.Item(14, 5).Value = "  -1.90%  "
$ws.Cells.Item(15, 2).Value = "WrappedBTC"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "68.469.86"
$ws.Cells.Item(15, 5).Value = "  -1.39%  "
$ws.Cells.Item(16, 2).Value = "ShibaInu"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.0000171"
$ws.Cells.Item(16, 5).Value = "  -3.10%  "
$ws.Cells.Item(17, 2).Value = "Avalanche"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "23.72"
$ws.Cells.Item(17, 5).Value = "  -4.54%  "
$ws.Cells.Item(18, 2).Value = "WrappedEther"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.433.92"
$ws.Cells.Item(18, 5).Value = "  -2.95%  "
$ws.Cells.Item(19, 2).Value = "Chainlink"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "10.92"
$ws.Cells.Item(19, 5).Value = "  -2.74%  "
$ws.Cells.Item(20, 2).Value = "BitcoinCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "352.50"
$ws.Cells.Item(20, 5).Value = "  +0.81%  "
$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "7.21"
$ws.Cells.Item(21, 5).Value = "  -4.56%  "
$ws.Cells.Item(22, 2).Value = "Polkadot"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.82"
$ws.Cells.Item(22, 5).Value = "  -2.34%  "
$ws.Cells.Item(23, 2).Value = "Dai"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "1.00"
$ws.Cells.Item(23, 5).Value = "  +0.00%  "
$ws.Cells.Item(24, 2).Value = "SuiNetwork"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.87"
$ws.Cells.Item(24, 5).Value = "  -4.39%  "
$ws.Cells.Item(25, 2).Value = "Litecoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "68.34"
$ws.Cells.Item(25, 5).Value = "  -2.72%  "
$ws.Cells.Item(26, 2).Value = "NEARProtocol"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.76"
$ws.Cells.Item(26, 5).Value = "  -4.70%  "
$ws.Cells.Item(27, 2).Value = "WrappedeETH"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.574.89"
$ws.Cells.Item(27, 5).Value = "  -2.49%  "
$ws.Cells.Item(28, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.02"
$ws.Cells.Item(28, 5).Value = "  +2.12%  "
$ws.Cells.Item(29, 2).Value = "Aptos"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.31"
$ws.Cells.Item(29, 5).Value = "  -5.97%  "
$ws.Cells.Item(30, 2).Value = "PEPE"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0₃0841"
$ws.Cells.Item(30, 5).Value = "  -5.34%  "
$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.34"
$ws.Cells.Item(31, 5).Value = "  -6.64%  "
$ws.Cells.Item(32, 5).Value = "  -3.41%  "
$ws.Cells.Item(33, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.00"
$ws.Cells.Item(33, 5).Value = "  +0.02%  "
$ws.Cells.Item(34, 2).Value = "Bittensor"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "430.93"
$ws.Cells.Item(34, 5).Value = "  -6.40%  "
$ws.Cells.Item(35, 2).Value = "PancakeSwap"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.69"
$ws.Cells.Item(35, 5).Value = "  -2.78%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "157.12"
$ws.Cells.Item(36, 5).Value = "  -1.71%  "
$ws.Cells.Item(37, 2).Value = "POPCAT"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.85"
$ws.Cells.Item(37, 5).Value = "  +92.88%  "
$ws.Cells.Item(38, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "19.01"
$ws.Cells.Item(38, 5).Value = "  -0.29%  "
$ws.Cells.Item(39, 5).Value = "  +0.04%  "
$ws.Cells.Item(40, 2).Value = "Kaspa"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.110"
$ws.Cells.Item(40, 5).Value = "  -5.68%  "
$ws.Cells.Item(41, 2).Value = "EthereumClassic"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "17.97"
$ws.Cells.Item(41, 5).Value = "  -2.58%  "
$ws.Cells.Item(42, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.307"
$ws.Cells.Item(42, 5).Value = "  -3.47%  "
$ws.Cells.Item(43, 2).Value = "RenderToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "4.52"
$ws.Cells.Item(43, 5).Value = "  -3.25%  "
$ws.Cells.Item(44, 2).Value = "Stacks"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.54"
$ws.Cells.Item(44, 5).Value = "  -3.69%  "
$ws.Cells.Item(45, 2).Value = "ImmutableX"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.09"
$ws.Cells.Item(45, 5).Value = "  +0.37%  "
$ws.Cells.Item(46, 2).Value = "dogwifhat"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.09"
$ws.Cells.Item(46, 5).Value = "  -5.74%  "
$ws.Cells.Item(47, 2).Value = "Aave"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "135.24"
$ws.Cells.Item(47, 5).Value = "  -5.01%  "
$ws.Cells.Item(48, 2).Value = "Filecoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.37"
$ws.Cells.Item(48, 5).Value = "  -2.77%  "
$ws.Cells.Item(49, 2).Value = "ARBITRUM"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.491"
$ws.Cells.Item(49, 5).Value = "  -5.34%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0717"
$ws.Cells.Item(50, 5).Value = "  -2.22%  "
$ws.Cells.Item(51, 2).Value = "Mantle"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.564"
$ws.Cells.Item(51, 5).Value = "  -2.19%  "
